$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the window slightly (cosmetic workbook view position change: xWindow 8100 -> 8000)
$win = $wb.Windows.Item(1)
$win.Left = 8000
$win.Top = 840

# Add new condition row 38: isSailFlagYes / 是否出海 / cache / sailFlag / != / number / 0
$ws.Cells.Item(38, 1).Value = "isSailFlagYes"
$ws.Cells.Item(38, 2).Value = "是否出海"
$ws.Cells.Item(38, 3).Value = "cache"
$ws.Cells.Item(38, 4).Value = "sailFlag"
$ws.Cells.Item(38, 5).Value = "!="
$ws.Cells.Item(38, 6).Value = "number"
$ws.Cells.Item(38, 7).Value = 0

# Update the active selection to C35
$ws.Range("C35").Select()
